$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Helper-ish approach: set a "text" value into a cell while forcing Excel to
# keep it as literal text (avoiding auto-conversion to number/currency for
# numeric-looking strings), then strip the temporary number-format style so
# the cell keeps the default style.
function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Data for the new rows 510-519 (appended to the bottom of Sheet1)
$rows = @(
    @{ r = 510; A = "21TRD09386"; B = "Bunner";   C = "DUS UCM";                                 D = "4510.111";     E = "UCM"; F = "No Contest"; G = "Guilty"; H = "$ 0";  I = "$ 0" },
    @{ r = 511; A = "21TRD09386"; B = "Bunner";   C = "TAIL LIGHTS-REAR LICENSE PLATE";          D = "4513.05";      E = "MM";  F = "No Contest"; G = "Guilty"; H = "$ 0";  I = "$ 0" },
    @{ r = 512; A = "21TRD09386"; B = "Bunner";   C = "DUS UCM";                                 D = "4510.111";     E = "UCM"; F = "No Contest"; G = "Guilty"; H = "$ 0";  I = "$ 0" },
    @{ r = 513; A = "21TRD09386"; B = "Bunner";   C = "TAIL LIGHTS-REAR LICENSE PLATE";          D = "4513.05";      E = "MM";  F = "No Contest"; G = "Guilty"; H = "$ 0";  I = "$ 0" },
    @{ r = 514; A = "21TRD09386"; B = "Bunner";   C = "DUS UCM";                                 D = "4510.111";     E = "UCM"; F = "No Contest"; G = "Guilty"; H = "$ 0";  I = "$ 0" },
    @{ r = 515; A = "21TRD09386"; B = "Bunner";   C = "TAIL LIGHTS-REAR LICENSE PLATE";          D = "4513.05";      E = "MM";  F = "No Contest"; G = "Guilty"; H = "$ 0";  I = "$ 0" },
    @{ r = 516; A = "21TRD09386"; B = "Bunner";   C = "DUS UCM";                                 D = "4510.111";     E = "UCM"; F = "No Contest"; G = "Guilty"; H = "$ 0";  I = "$ 0" },
    @{ r = 517; A = "21TRD09386"; B = "Bunner";   C = "TAIL LIGHTS-REAR LICENSE PLATE";          D = "4513.05";      E = "MM";  F = "No Contest"; G = "Guilty"; H = "$ 0";  I = "$ 0" },
    @{ r = 518; A = "21CRB01291"; B = "Bunner";   C = "PERMISSION REQ'D TO USE LICENSED DOCK";   D = "1501:46-12-04";E = "MM";  F = "Guilty";     G = "Guilty"; H = "$ 0";  I = "$ 0" },
    @{ r = 519; A = "21CRB01268"; B = "Hemmeter"; C = "POSSESSION DRUG PARAPHERNALIA";           D = "2925.14(C)";   E = "M4";  F = "Guilty";     G = "Guilty"; H = "$ 50"; I = "$ 25"; J = "10"; K = "5" }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    Set-TextValue $r 4 $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    Set-TextValue $r 8 $row.H
    Set-TextValue $r 9 $row.I
    if ($row.ContainsKey("J")) {
        Set-TextValue $r 10 $row.J
    }
    if ($row.ContainsKey("K")) {
        Set-TextValue $r 11 $row.K
    }
}
